$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the entire "21 Preparing for Cancer Treatment" section: the
#    Heading2 paragraph plus its four one-line bullet paragraphs
#    (Primary Care Physician / MyAtrium Portal / Exercise / Smoking
#    Cessation). This whole block is removed by the target edit.
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*21 Preparing for Cancer Treatment*") {
        $startPara = $i
    }
    if ($startPara -ne $null -and $t -like "*Smoking Cessation*" -and $i -gt $startPara -and $endPara -eq $null) {
        $endPara = $i
    }
    if ($startPara -ne $null -and $endPara -ne $null) {
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rStart = $d.Paragraphs.Item($startPara).Range.Start
    $rEnd = $d.Paragraphs.Item($endPara).Range.End
    $killRange = $d.Range($rStart, $rEnd)
    $killRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Renumber the following four Heading2 captions down by one, now that the
#    "Preparing for Cancer Treatment" heading is gone.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("22 Primary Care Practitioner (PCP)", $true, $true, $false, $false, $false, $true, 1, $false, "21 Primary Care Practitioner (PCP)", 2) | Out-Null
$d.Content.Find.Execute("23 My Atrium Patient Portal", $true, $true, $false, $false, $false, $true, 1, $false, "22 My Atrium Patient Portal", 2) | Out-Null
$d.Content.Find.Execute("24 Exercise", $true, $true, $false, $false, $false, $true, 1, $false, "23 Exercise", 2) | Out-Null
$d.Content.Find.Execute("25 Smoking Cessation", $true, $true, $false, $false, $false, $true, 1, $false, "24 Smoking Cessation", 2) | Out-Null

Write-Host "edit complete"
